$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(705.698026459629,    1157.793623994623,   1444.187648963287,   1589.319899554724),
    @(713.4871374663375,   1167.018325727053,   1450.714238129294,   1593.901005905419),
    @(636.8114506417281,   1102.723238839376,   1409.905763062105,   1555.648480433783),
    @(714.0936209589279,   1176.781486762983,   1454.454665392003,   1604.157106965235),
    @(671.7520134723012,   1140.941785076102,   1438.279270177918,   1592.155716468016),
    @(750.8099286789158,   1207.976052222934,   1483.315502256271,   1605.734104783698),
    @(356.6613353955055,   807.3842191477663,   1181.738694825472,   1465.874424671512),
    @(665.5488374974055,   1118.584399760819,   1420.839406625197,   1566.306438240137),
    @(796.7028593535204,   1256.002499623704,   1513.159627900369,   1641.04387868246),
    @(788.236008520571,    1248.948808275188,   1508.582765334469,   1637.111851400533),
    @(720.6201348792119,   1192.728605902402,   1479.979521199734,   1605.589635760548),
    @(774.2730251325225,   1237.533746002944,   1501.870447865611,   1631.387590119088)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rowValues[0]
    $ws.Cells.Item($row, 3).Value = $rowValues[1]
    $ws.Cells.Item($row, 4).Value = $rowValues[2]
    $ws.Cells.Item($row, 5).Value = $rowValues[3]
}
